# Apply updated cryptocurrency price/volume data to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a literal text value into a cell without letting Excel
# auto-convert numeric-looking strings (e.g. "233.82") into real numbers.
function Set-TextValue {
    param($rangeAddr, $val)
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "37.339.28"
$ws.Range("E2").Value = "  -0.62%  "
Set-TextValue "D3" "2.066.89"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue "D5" "233.82"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  +0.82%  "
Set-TextValue "D8" "56.72"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("E9").Value = "  -0.79%  "
Set-TextValue "D10" "0.0762"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +0.39%  "
Set-TextValue "D12" "2.369.47"
$ws.Range("E12").Value = "  -0.54%  "
Set-TextValue "D13" "14.59"
$ws.Range("E13").Value = "  +0.06%  "
Set-TextValue "D14" "20.62"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("E15").Value = "  -0.57%  "
Set-TextValue "D16" "5.11"
$ws.Range("E16").Value = "  -2.86%  "
Set-TextValue "D17" "2.063.57"
$ws.Range("E17").Value = "  -3.97%  "
Set-TextValue "D18" "37.264.49"
$ws.Range("E18").Value = "  -1.08%  "
Set-TextValue "D19" "6.27"
$ws.Range("E19").Value = "  +4.01%  "
Set-TextValue "D20" "69.49"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  -0.54%  "
Set-TextValue "D22" "226.18"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  -2.30%  "
Set-TextValue "D26" "167.47"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +2.86%  "
Set-TextValue "D29" "19.01"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  -1.18%  "
Set-TextValue "D33" "0.0615"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  +3.36%  "
Set-TextValue "D35" "2.49"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -4.39%  "
Set-TextValue "D39" "5.64"
$ws.Range("E39").Value = "  -5.34%  "
$ws.Range("E40").Value = "  -0.21%  "
Set-TextValue "D41" "1.469.53"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D42" "0.0939"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "96.06"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D44" "4.32"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  +2.92%  "
Set-TextValue "D46" "0.0212"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  -1.46%  "
Set-TextValue "D48" "14.98"
$ws.Range("E48").Value = "  -9.46%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  -3.03%  "
Set-TextValue "D51" "2.258.13"
$ws.Range("E51").Value = "  -0.48%  "
